# Scheduled runner update: refresh computed market-price / profit figures
# for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1136.3334
$ws.Range("J80").Value = 1252.5555
$ws.Range("L80").Value = 3757.6665
$ws.Range("N80").Value = -5753.666499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1136.3334
$ws.Range("J83").Value = 1252.5555
$ws.Range("L83").Value = 11272.9995
$ws.Range("N83").Value = -21256.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 992.8946999999999
$ws.Range("I137").Value = 689.8
$ws.Range("K137").Value = 2069.4
$ws.Range("M137").Value = 480.6000000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4971.5557
$ws.Range("I138").Value = 3807.5454
$ws.Range("J138").Value = 5771.8125
$ws.Range("K138").Value = 11422.6362
$ws.Range("L138").Value = 17315.4375
$ws.Range("M138").Value = -6282.636200000001
$ws.Range("N138").Value = -27595.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3875
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 4500
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 4500
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -4726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9269.6
$ws.Range("I45").Value = 9269.6
$ws.Range("K45").Value = 9269.6
$ws.Range("M45").Value = -8892.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 958.0833
$ws.Range("I74").Value = 961.2857
$ws.Range("K74").Value = 961.2857
$ws.Range("M74").Value = -87.28570000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 958.0833
$ws.Range("I77").Value = 961.2857
$ws.Range("K77").Value = 4806.4285
$ws.Range("M77").Value = -438.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3875
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = 294
$ws.Range("N116").Value = -9088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3875
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = -1886
$ws.Range("N3").Value = -4728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3603.3333
$ws.Range("I99").Value = 3655
$ws.Range("K99").Value = 3655
$ws.Range("M99").Value = -2157

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1654.4722
$ws.Range("I31").Value = 1096.3334
$ws.Range("J31").Value = 1933.5416
$ws.Range("K31").Value = 1096.3334
$ws.Range("L31").Value = 1933.5416
$ws.Range("M31").Value = -801.3334
$ws.Range("N31").Value = -2523.5416

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1654.4722
$ws.Range("I34").Value = 1096.3334
$ws.Range("J34").Value = 1933.5416
$ws.Range("K34").Value = 1096.3334
$ws.Range("L34").Value = 1933.5416
$ws.Range("M34").Value = -894.3334
$ws.Range("N34").Value = -2337.5416

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3518
$ws.Range("I58").Value = 2870.2
$ws.Range("K58").Value = 2870.2
$ws.Range("M58").Value = -2667.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4949
$ws.Range("I86").Value = 4949
$ws.Range("K86").Value = 4949
$ws.Range("M86").Value = -3826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4949
$ws.Range("I89").Value = 4949
$ws.Range("K89").Value = 24745
$ws.Range("M89").Value = -19129

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 62711.25
$ws.Range("J108").Value = 62711.25
$ws.Range("L108").Value = 62711.25
$ws.Range("N108").Value = -70391.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4188
$ws.Range("J122").Value = 5282.5
$ws.Range("L122").Value = 15847.5
$ws.Range("N122").Value = -20747.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4654.4546
$ws.Range("I132").Value = 4654.4546
$ws.Range("K132").Value = 13963.3638
$ws.Range("M132").Value = -11433.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2065.8333
$ws.Range("I134").Value = 2036.5625
$ws.Range("K134").Value = 6109.6875
$ws.Range("M134").Value = -3574.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3518
$ws.Range("I136").Value = 2870.2
$ws.Range("K136").Value = 8610.599999999999
$ws.Range("M136").Value = -6060.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 293.75
$ws.Range("I99").Value = 293.75
$ws.Range("K99").Value = 881.25
$ws.Range("M99").Value = 1364.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2752501
$ws.Range("J11").Value = 502502
$ws.Range("L11").Value = 502502
$ws.Range("N11").Value = -502780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 167500
$ws.Range("I46").Value = 275000
$ws.Range("K46").Value = 275000
$ws.Range("M46").Value = -274844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998
$ws.Range("I70").Value = 4997
$ws.Range("K70").Value = 4997
$ws.Range("M70").Value = -4727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4998
$ws.Range("I73").Value = 4997
$ws.Range("K73").Value = 4997
$ws.Range("M73").Value = -4061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25050.5
$ws.Range("J80").Value = 29760.6
$ws.Range("L80").Value = 29760.6
$ws.Range("N80").Value = -31756.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 25050.5
$ws.Range("J83").Value = 29760.6
$ws.Range("L83").Value = 148803
$ws.Range("N83").Value = -158787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5800.1
$ws.Range("I126").Value = 2745.5
$ws.Range("K126").Value = 8236.5
$ws.Range("M126").Value = -5766.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5045.3335
$ws.Range("J132").Value = 1500.5
$ws.Range("L132").Value = 4501.5
$ws.Range("N132").Value = -9561.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2479.75
$ws.Range("I22").Value = 1730
$ws.Range("J22").Value = 3729.3333
$ws.Range("K22").Value = 1730
$ws.Range("L22").Value = 3729.3333
$ws.Range("M22").Value = -1435
$ws.Range("N22").Value = -4319.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2479.75
$ws.Range("I27").Value = 1730
$ws.Range("J27").Value = 3729.3333
$ws.Range("K27").Value = 1730
$ws.Range("L27").Value = 3729.3333
$ws.Range("M27").Value = -1623
$ws.Range("N27").Value = -3943.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1696.3334
$ws.Range("I46").Value = 1957
$ws.Range("J46").Value = 1596.0769
$ws.Range("K46").Value = 1957
$ws.Range("L46").Value = 1596.0769
$ws.Range("M46").Value = -1769
$ws.Range("N46").Value = -1972.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2445.1304
$ws.Range("I132").Value = 2105.6667
$ws.Range("K132").Value = 6317.000100000001
$ws.Range("M132").Value = -3787.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2658.6
$ws.Range("J136").Value = 2846.5
$ws.Range("L136").Value = 8539.5
$ws.Range("N136").Value = -13639.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1399.3334
$ws.Range("I132").Value = 1424.5
$ws.Range("J132").Value = 1349
$ws.Range("K132").Value = 4273.5
$ws.Range("L132").Value = 4047
$ws.Range("M132").Value = -1743.5
$ws.Range("N132").Value = -9107

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2767.6
$ws.Range("I136").Value = 2630.6667
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 7892.000100000001
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -5342.000100000001
$ws.Range("N136").Value = -17100
